# Commit: "skip empty lines" + "add empty columns even if no data is present"
#
# The "jack" row (previously row 3) becomes blank/skipped, so the data that
# followed (the "helen" row, previously row 4) shifts down by one. We get
# this shape by inserting a blank row above the old row 3 (pushing jack -> 4
# and helen -> 5) and then clearing that newly inserted row so it carries no
# cells/formatting and disappears from the saved sheetData entirely.
#
# A new header column "emptyColumn" is also introduced at D1, even though
# the column itself has no data rows underneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the "jack" row, shifting jack/helen down one row,
# then clear it completely so the "skipped" row leaves no trace.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Clear()

# New header for the empty column (D) - no data underneath it.
$ws.Range("D1").Value = "emptyColumn"

# Match the author's final selection.
$ws.Range("B12").Select()
